$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new header labels
$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

# Copy the formatting from an existing header cell (AE1) so the new headers match
$ws.Range("AE1").Copy() | Out-Null
$ws.Range("AF1:AH1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($r = 2; $r -le 55; $r++) {
    $ypredH = $ws.Cells.Item($r, 23).Value2  # W
    $ypredD = $ws.Cells.Item($r, 24).Value2  # X
    $ypredA = $ws.Cells.Item($r, 25).Value2  # Y
    $ytrueH = $ws.Cells.Item($r, 26).Value2  # Z
    $ytrueD = $ws.Cells.Item($r, 27).Value2  # AA
    $ytrueA = $ws.Cells.Item($r, 28).Value2  # AB

    $pctDiffH = (($ypredH - $ytrueH) / $ytrueH) * 100
    $pctDiffD = (($ypredD - $ytrueD) / $ytrueD) * 100
    $pctDiffA = (($ypredA - $ytrueA) / $ytrueA) * 100

    $ws.Cells.Item($r, 32).Value = $pctDiffH  # AF
    $ws.Cells.Item($r, 33).Value = $pctDiffD  # AG
    $ws.Cells.Item($r, 34).Value = $pctDiffA  # AH
}
